$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.423.09'
$ws.Range("E2").Value = '  -7.25%  '
$ws.Range("D3").Value = '1.439.30'
$ws.Range("E3").Value = '  -7.35%  '
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '277.65'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3752'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.54%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3062'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '40.31'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -8.34%  '
$ws.Range("E10").Value = '  -3.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06539'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -8.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.004'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.343'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.23'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -6.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.113'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -6.99%  '
$ws.Range("D16").Value = '1.441.70'
$ws.Range("E16").Value = '  -6.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001008'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -7.51%  '
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '76.21'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -7.54%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05855'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -10.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.712'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -6.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.34'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -5.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.80'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.314'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.50%  '
$ws.Range("D25").Value = '20.417.83'
$ws.Range("E25").Value = '  -7.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.83'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.193'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -6.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.95'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -6.85%  '
$ws.Range("D29").Value = '1.604.03'
$ws.Range("E29").Value = '  -7.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '109.55'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.839'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -21.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8972'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -6.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.400'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -6.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07712'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -6.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.309'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -7.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.003'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.77'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +3.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05649'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.139'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.702'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -6.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02039'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -8.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1912'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.336'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -16.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.581'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5303'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -6.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.16'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5135'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -6.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '112.40'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.784'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.052'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -5.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.004'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.04%  '
